$d = $word.ActiveDocument

# Replace all occurrences of "June 07, 2022" with "June 08, 2022"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("June 07, 2022", $true, $true, $false, $false, $false, $true, 1, $false, "June 08, 2022", 2)

# Replace the single occurrence of "August 06, 2022" with "August 07, 2022"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("August 06, 2022", $true, $true, $false, $false, $false, $true, 1, $false, "August 07, 2022", 2)
